# Translate the customer & advisor login-entry strings, and add new
# key/value rows for the login page labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the previously-untranslated nav.* rows (B column) ---
$ws.Cells.Item(102, 2).Value = "קצת עלינו"
$ws.Cells.Item(103, 2).Value = "שאלות נפוצות"
$ws.Cells.Item(104, 2).Value = "מחירון"
$ws.Cells.Item(105, 2).Value = "אבטחה ופרטיות"
$ws.Cells.Item(106, 2).Value = "כניסת לקוחות"
$ws.Cells.Item(107, 2).Value = "כניסת יועצים/סוכנים"

# --- Append new key/value/format_flag rows for the login page ---
$ws.Cells.Item(108, 1).Value = "login.label"
$ws.Cells.Item(108, 2).Value = "!בוש ךתוארל םיחמש"
$ws.Cells.Item(108, 3).Value = 1

$ws.Cells.Item(109, 1).Value = "login.entry_customer"
$ws.Cells.Item(109, 2).Value = "תוחוקל תסינכ"
$ws.Cells.Item(109, 3).Value = 1

$ws.Cells.Item(110, 1).Value = "login.entry_advisor"
$ws.Cells.Item(110, 2).Value = "םינכוס/םיצעוי תסינכ"
$ws.Cells.Item(110, 3).Value = 1

$ws.Cells.Item(111, 1).Value = "login.id_card_number.label"
$ws.Cells.Item(111, 2).Value = "תוהז תדועת"
$ws.Cells.Item(111, 3).Value = 1

$ws.Cells.Item(112, 1).Value = "login.phone.label"
$ws.Cells.Item(112, 2).Value = "ןופלט רפסמ"
$ws.Cells.Item(112, 3).Value = 1

# --- Match the author's final view state (scroll + selection) ---
[void]$ws.Range("L23").Select()
